$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "womens capri leggings xl"
$ws.Range("A2").Value = "womens capri overalls"
$ws.Range("A3").Value = "womens capri pants"
$ws.Range("A4").Value = "womens capri pants size 12"
$ws.Range("A5").Value = "womens capri running tights"
$ws.Range("A6").Value = "womens capri shorts"
$ws.Range("A7").Value = "womens capri tights"
$ws.Range("A8").Value = "womens capri workout pants"
$ws.Range("A9").Value = "womens capri yoga pants"
$ws.Range("A10").Value = "womens capri yoga workout running gym pants"
$ws.Range("A11").Value = "womens capris"
$ws.Range("A12").Value = "womens capris leggings"
$ws.Range("A13").Value = "womens capris with pockets"
$ws.Range("A14").Value = "womens catching gear"
$ws.Range("A15").Value = "womens chillys"
$ws.Range("A16").Value = "womens christmas tights"
$ws.Range("A17").Value = "womens clothing tape"
$ws.Range("A18").Value = "womens cold compression tops"
$ws.Range("A19").Value = "womens cold gear compression leggings"
$ws.Range("A20").Value = "womens cold running gear"
$ws.Range("A21").Value = "womens cold weather biking gear"
$ws.Range("A22").Value = "womens cold weather dresses"
$ws.Range("A23").Value = "womens cold weather leggings"
$ws.Range("A24").Value = "womens cold weather running gear"
$ws.Range("A25").Value = "womens cold weather running tops"
$ws.Range("A26").Value = "womens cold weather workout gear"
$ws.Range("A27").Value = "womens compresion shorts"
$ws.Range("A28").Value = "womens compression"
$ws.Range("A29").Value = "womens compression capri leggings"
$ws.Range("A30").Value = "womens compression capri pants"
$ws.Range("A31").Value = "womens compression capris"
$ws.Range("A32").Value = "womens compression clothing"
$ws.Range("A33").Value = "womens compression knee highs"
$ws.Range("A34").Value = "womens compression legging"
$ws.Range("A35").Value = "womens compression leggings"
$ws.Range("A36").Value = "womens compression leggings high waist"
$ws.Range("A37").Value = "womens compression leggings plus size"
$ws.Range("A38").Value = "womens compression leggings white"
$ws.Range("A39").Value = "womens compression pant"
$ws.Range("A40").Value = "womens compression pants"
$ws.Range("A41").Value = "women's compression pants & tights"
$ws.Range("A42").Value = "womens compression pants for running"
$ws.Range("A43").Value = "womens compression pants high waist"
$ws.Range("A44").Value = "womens compression pants plus size"
$ws.Range("A45").Value = "womens compression pants white"
$ws.Range("A46").Value = "womens compression running capris"
$ws.Range("A47").Value = "womens compression running leggings"
$ws.Range("A48").Value = "womens compression running pants"
$ws.Range("A49").Value = "womens compression running short"
$ws.Range("A50").Value = "womens compression running shorts"
$ws.Range("A51").Value = "womens compression running tights"
$ws.Range("A52").Value = "womens compression short"
$ws.Range("A53").Value = "womens compression shorts"
$ws.Range("A54").Value = "womens compression shorts 3 inch"
$ws.Range("A55").Value = "womens compression shorts 6 inch"
$ws.Range("A56").Value = "womens compression shorts 7 inch"
$ws.Range("A57").Value = "womens compression shorts for running"
$ws.Range("A58").Value = "womens compression shorts high waist"
$ws.Range("A59").Value = "womens compression shorts long"
$ws.Range("A60").Value = "womens compression shorts nike"
$ws.Range("A61").Value = "womens compression suit"
$ws.Range("A62").Value = "womens compression swim capris"
$ws.Range("A63").Value = "womens compression tank"
$ws.Range("A64").Value = "womens compression tank top long"
$ws.Range("A65").Value = "womens compression tanks"
$ws.Range("A66").Value = "womens compression tights"
$ws.Range("A67").Value = "womens compression tights leggings"
$ws.Range("A68").Value = "womens compression tops"
$ws.Range("A69").Value = "womens compression underwear"
$ws.Range("A70").Value = "womens compression wear"
$ws.Range("A71").Value = "womens compression workout"
$ws.Range("A72").Value = "womens compression workout clothes"
$ws.Range("A73").Value = "womens compression workout leggings"
$ws.Range("A74").Value = "womens compression workout pants"
$ws.Range("A75").Value = "womens compression yoga"
$ws.Range("A76").Value = "womens compression yoga pants"
$ws.Range("A77").Value = "womens compressionz"
$ws.Range("A78").Value = "womens conpression shorts"
$ws.Range("A79").Value = "womens cothes"
$ws.Range("A80").Value = "womens cwx"
$ws.Range("A81").Value = "womens cycling capris"
$ws.Range("A82").Value = "womens cycling short"
$ws.Range("A83").Value = "womens cycling shorts"
$ws.Range("A84").Value = "womens cycling tights"
$ws.Range("A85").Value = "womens double layer rainbow"
$ws.Range("A86").Value = "womens drawstring shorts"
$ws.Range("A87").Value = "womens dress tights"
$ws.Range("A88").Value = "womens elastic waist pants size 14"
$ws.Range("A89").Value = "womens exercise apparel"
$ws.Range("A90").Value = "womens exercise capri pants"
$ws.Range("A91").Value = "womens exercise clothes"
$ws.Range("A92").Value = "womens exercise pants"
$ws.Range("A93").Value = "womens exercise wear"
$ws.Range("A94").Value = "womens fashion tights"
$ws.Range("A95").Value = "womens fitness bike"
$ws.Range("A96").Value = "womens fitness clothes"
$ws.Range("A97").Value = "womens fitness clothing"
$ws.Range("A98").Value = "womens fitted leggings"
$ws.Range("A99").Value = "womens football pants"
$ws.Range("A100").Value = "womens free run distance"
